# Apply the "merge Intel-Irris recent state" data updates to the watermark
# workbook: update the measured resistance (D2) and the third calibration
# temperature (A7). Dependent formulas in D5:F7 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1044
$ws.Range("A7").Value = 9

$excel.Calculate()
